$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 37.17328633333333
$ws.Cells.Item(2, 8).Value = 111.519859
$ws.Cells.Item(2, 9).Value = 0.005170079968594893
$ws.Cells.Item(2, 10).Value = 0.005188590814393131
$ws.Cells.Item(2, 13).Value = 10.250695
$ws.Cells.Item(2, 14).Value = 30.752085
$ws.Cells.Item(2, 15).Value = 0.2501330709220828
$ws.Cells.Item(2, 16).Value = 0.259830908271274
$ws.Cells.Item(2, 17).Value = 381.0520203506683
$ws.Cells.Item(2, 18).Value = 3429.468183156015
$ws.Cells.Item(2, 19).Value = 0.001293207979457386
$ws.Cells.Item(2, 20).Value = 0.001348156263951756

$ws.Cells.Item(3, 7).Value = 37.17328633333333
$ws.Cells.Item(3, 8).Value = 111.519859
$ws.Cells.Item(3, 9).Value = 0.005170079968594893
$ws.Cells.Item(3, 10).Value = 0.005188590814393131
$ws.Cells.Item(3, 15).Value = 0.2460808482884365
$ws.Cells.Item(3, 16).Value = 0.2556215780794031
$ws.Cells.Item(3, 17).Value = 374.8788757290096
$ws.Cells.Item(3, 18).Value = 3373.909881561087
$ws.Cells.Item(3, 19).Value = 0.001272257664390884
$ws.Cells.Item(3, 20).Value = 0.001326315771983467

$ws.Cells.Item(4, 7).Value = 37.17328633333333
$ws.Cells.Item(4, 8).Value = 111.519859
$ws.Cells.Item(4, 9).Value = 0.005170079968594893
$ws.Cells.Item(4, 10).Value = 0.005188590814393131
$ws.Cells.Item(4, 13).Value = 7.311799000000001
$ws.Cells.Item(4, 14).Value = 21.935397
$ws.Cells.Item(4, 15).Value = 0.1784193889131434
$ws.Cells.Item(4, 16).Value = 0.1853368357235283
$ws.Cells.Item(4, 17).Value = 271.8035978387803
$ws.Cells.Item(4, 18).Value = 2446.232380549023
$ws.Cells.Item(4, 19).Value = 0.0009224425086287843
$ws.Cells.Item(4, 20).Value = 0.0009616370034037878

$ws.Cells.Item(5, 7).Value = 37.17328633333333
$ws.Cells.Item(5, 8).Value = 111.519859
$ws.Cells.Item(5, 9).Value = 0.005170079968594893
$ws.Cells.Item(5, 10).Value = 0.005188590814393131
$ws.Cells.Item(5, 13).Value = 4.588677499999999
$ws.Cells.Item(5, 14).Value = 9.177354999999999
$ws.Cells.Item(5, 15).Value = 0.1119709438770801
$ws.Cells.Item(5, 16).Value = 0.07754142475796089
$ws.Cells.Item(5, 17).Value = 170.5762225988241
$ws.Cells.Item(5, 18).Value = 1023.457335592945
$ws.Cells.Item(5, 19).Value = 0.0005788987340035548
$ws.Cells.Item(5, 20).Value = 0.000402330724234112

$ws.Cells.Item(6, 7).Value = 37.17328633333333
$ws.Cells.Item(6, 8).Value = 111.519859
$ws.Cells.Item(6, 9).Value = 0.005170079968594893
$ws.Cells.Item(6, 10).Value = 0.005188590814393131
$ws.Cells.Item(6, 13).Value = 8.745164000000001
$ws.Cells.Item(6, 14).Value = 26.235492
$ws.Cells.Item(6, 15).Value = 0.2133957479992572
$ws.Cells.Item(6, 16).Value = 0.2216692531678338
$ws.Cells.Item(6, 17).Value = 325.0864854039586
$ws.Cells.Item(6, 18).Value = 2925.778368635628
$ws.Cells.Item(6, 19).Value = 0.001103273082114283
$ws.Cells.Item(6, 20).Value = 0.001150151050820008

$ws.Cells.Item(7, 9).Value = 0.006280726092526873
$ws.Cells.Item(7, 10).Value = 0.006303213472394487
$ws.Cells.Item(7, 13).Value = 10.250695
$ws.Cells.Item(7, 14).Value = 30.752085
$ws.Cells.Item(7, 15).Value = 0.2501330709220828
$ws.Cells.Item(7, 16).Value = 0.259830908271274
$ws.Cells.Item(7, 17).Value = 462.9103188662983
$ws.Cells.Item(7, 18).Value = 4166.192869796685
$ws.Cells.Item(7, 19).Value = 0.0015710173051442
$ws.Cells.Item(7, 20).Value = 0.00163776968155999

$ws.Cells.Item(8, 9).Value = 0.006280726092526873
$ws.Cells.Item(8, 10).Value = 0.006303213472394487
$ws.Cells.Item(8, 15).Value = 0.2460808482884365
$ws.Cells.Item(8, 16).Value = 0.2556215780794031
$ws.Cells.Item(8, 19).Value = 0.00154556640471633
$ws.Cells.Item(8, 20).Value = 0.001611237374784833

$ws.Cells.Item(9, 9).Value = 0.006280726092526873
$ws.Cells.Item(9, 10).Value = 0.006303213472394487
$ws.Cells.Item(9, 13).Value = 7.311799000000001
$ws.Cells.Item(9, 14).Value = 21.935397
$ws.Cells.Item(9, 15).Value = 0.1784193889131434
$ws.Cells.Item(9, 16).Value = 0.1853368357235283
$ws.Cells.Item(9, 17).Value = 330.1929485343464
$ws.Cells.Item(9, 18).Value = 2971.736536809117
$ws.Cells.Item(9, 19).Value = 0.001120603311359479
$ws.Cells.Item(9, 20).Value = 0.001168217639863508

$ws.Cells.Item(10, 9).Value = 0.006280726092526873
$ws.Cells.Item(10, 10).Value = 0.006303213472394487
$ws.Cells.Item(10, 13).Value = 4.588677499999999
$ws.Cells.Item(10, 14).Value = 9.177354999999999
$ws.Cells.Item(10, 15).Value = 0.1119709438770801
$ws.Cells.Item(10, 16).Value = 0.07754142475796089
$ws.Cells.Item(10, 17).Value = 207.2197216578591
$ws.Cells.Item(10, 18).Value = 1243.318329947155
$ws.Cells.Item(10, 19).Value = 0.000703258828813639
$ws.Cells.Item(10, 20).Value = 0.0004887601532030426

$ws.Cells.Item(11, 9).Value = 0.006280726092526873
$ws.Cells.Item(11, 10).Value = 0.006303213472394487
$ws.Cells.Item(11, 13).Value = 8.745164000000001
$ws.Cells.Item(11, 14).Value = 26.235492
$ws.Cells.Item(11, 15).Value = 0.2133957479992572
$ws.Cells.Item(11, 16).Value = 0.2216692531678338
$ws.Cells.Item(11, 17).Value = 394.9221643779347
$ws.Cells.Item(11, 18).Value = 3554.299479401413
$ws.Cells.Item(11, 19).Value = 0.001340280242493224
$ws.Cells.Item(11, 20).Value = 0.001397228622983114

$ws.Cells.Item(12, 7).Value = 2375.59786
$ws.Cells.Item(12, 8).Value = 7126.79358
$ws.Cells.Item(12, 9).Value = 0.3303993840977568
$ws.Cells.Item(12, 10).Value = 0.3315823391174117
$ws.Cells.Item(12, 13).Value = 10.250695
$ws.Cells.Item(12, 14).Value = 30.752085
$ws.Cells.Item(12, 15).Value = 0.2501330709220828
$ws.Cells.Item(12, 16).Value = 0.259830908271274
$ws.Cells.Item(12, 17).Value = 24351.52910551269
$ws.Cells.Item(12, 18).Value = 219163.7619496143
$ws.Cells.Item(12, 19).Value = 0.08264381257513667
$ws.Cells.Item(12, 20).Value = 0.08615534033959066

$ws.Cells.Item(13, 7).Value = 2375.59786
$ws.Cells.Item(13, 8).Value = 7126.79358
$ws.Cells.Item(13, 9).Value = 0.3303993840977568
$ws.Cells.Item(13, 10).Value = 0.3315823391174117
$ws.Cells.Item(13, 15).Value = 0.2460808482884365
$ws.Cells.Item(13, 16).Value = 0.2556215780794031
$ws.Cells.Item(13, 17).Value = 23957.02782248966
$ws.Cells.Item(13, 18).Value = 215613.2504024069
$ws.Cells.Item(13, 19).Value = 0.08130496071275296
$ws.Cells.Item(13, 20).Value = 0.08475960078845256

$ws.Cells.Item(14, 7).Value = 2375.59786
$ws.Cells.Item(14, 8).Value = 7126.79358
$ws.Cells.Item(14, 9).Value = 0.3303993840977568
$ws.Cells.Item(14, 10).Value = 0.3315823391174117
$ws.Cells.Item(14, 13).Value = 7.311799000000001
$ws.Cells.Item(14, 14).Value = 21.935397
$ws.Cells.Item(14, 15).Value = 0.1784193889131434
$ws.Cells.Item(14, 16).Value = 0.1853368357235283
$ws.Cells.Item(14, 17).Value = 17369.89405715014
$ws.Cells.Item(14, 18).Value = 156329.0465143513
$ws.Cells.Item(14, 19).Value = 0.05894965620800072
$ws.Cells.Item(14, 20).Value = 0.06145442151382699

$ws.Cells.Item(15, 7).Value = 2375.59786
$ws.Cells.Item(15, 8).Value = 7126.79358
$ws.Cells.Item(15, 9).Value = 0.3303993840977568
$ws.Cells.Item(15, 10).Value = 0.3315823391174117
$ws.Cells.Item(15, 13).Value = 4.588677499999999
$ws.Cells.Item(15, 14).Value = 9.177354999999999
$ws.Cells.Item(15, 15).Value = 0.1119709438770801
$ws.Cells.Item(15, 16).Value = 0.07754142475796089
$ws.Cells.Item(15, 17).Value = 10900.85244923015
$ws.Cells.Item(15, 18).Value = 65405.11469538089
$ws.Cells.Item(15, 19).Value = 0.03699513089383176
$ws.Cells.Item(15, 20).Value = 0.02571136699974145

$ws.Cells.Item(16, 7).Value = 2375.59786
$ws.Cells.Item(16, 8).Value = 7126.79358
$ws.Cells.Item(16, 9).Value = 0.3303993840977568
$ws.Cells.Item(16, 10).Value = 0.3315823391174117
$ws.Cells.Item(16, 13).Value = 8.745164000000001
$ws.Cells.Item(16, 14).Value = 26.235492
$ws.Cells.Item(16, 15).Value = 0.2133957479992572
$ws.Cells.Item(16, 16).Value = 0.2216692531678338
$ws.Cells.Item(16, 17).Value = 20774.99288374904
$ws.Cells.Item(16, 18).Value = 186974.9359537413
$ws.Cells.Item(16, 19).Value = 0.0705058237080347
$ws.Cells.Item(16, 20).Value = 0.07350160947580005

$ws.Cells.Item(17, 7).Value = 76.954105
$ws.Cells.Item(17, 8).Value = 153.90821
$ws.Cells.Item(17, 9).Value = 0.01070281688829022
$ws.Cells.Item(17, 10).Value = 0.007160758019481436
$ws.Cells.Item(17, 13).Value = 10.250695
$ws.Cells.Item(17, 14).Value = 30.752085
$ws.Cells.Item(17, 15).Value = 0.2501330709220828
$ws.Cells.Item(17, 16).Value = 0.259830908271274
$ws.Cells.Item(17, 17).Value = 788.8330593529748
$ws.Cells.Item(17, 18).Value = 4732.99835611785
$ws.Cells.Item(17, 19).Value = 0.002677128455784763
$ws.Cells.Item(17, 20).Value = 0.001860586260112671

$ws.Cells.Item(18, 7).Value = 76.954105
$ws.Cells.Item(18, 8).Value = 153.90821
$ws.Cells.Item(18, 9).Value = 0.01070281688829022
$ws.Cells.Item(18, 10).Value = 0.007160758019481436
$ws.Cells.Item(18, 15).Value = 0.2460808482884365
$ws.Cells.Item(18, 16).Value = 0.2556215780794031
$ws.Cells.Item(18, 17).Value = 776.053752860255
$ws.Cells.Item(18, 18).Value = 4656.32251716153
$ws.Cells.Item(18, 19).Value = 0.002633758258946262
$ws.Cells.Item(18, 20).Value = 0.001830444265184586

$ws.Cells.Item(19, 7).Value = 76.954105
$ws.Cells.Item(19, 8).Value = 153.90821
$ws.Cells.Item(19, 9).Value = 0.01070281688829022
$ws.Cells.Item(19, 10).Value = 0.007160758019481436
$ws.Cells.Item(19, 13).Value = 7.311799000000001
$ws.Cells.Item(19, 14).Value = 21.935397
$ws.Cells.Item(19, 15).Value = 0.1784193889131434
$ws.Cells.Item(19, 16).Value = 0.1853368357235283
$ws.Cells.Item(19, 17).Value = 562.672947984895
$ws.Cells.Item(19, 18).Value = 3376.03768790937
$ws.Cells.Item(19, 19).Value = 0.001909590048858012
$ws.Cells.Item(19, 20).Value = 0.001327152232712569

$ws.Cells.Item(20, 7).Value = 76.954105
$ws.Cells.Item(20, 8).Value = 153.90821
$ws.Cells.Item(20, 9).Value = 0.01070281688829022
$ws.Cells.Item(20, 10).Value = 0.007160758019481436
$ws.Cells.Item(20, 13).Value = 4.588677499999999
$ws.Cells.Item(20, 14).Value = 9.177354999999999
$ws.Cells.Item(20, 15).Value = 0.1119709438770801
$ws.Cells.Item(20, 16).Value = 0.07754142475796089
$ws.Cells.Item(20, 17).Value = 353.1175701461374
$ws.Cells.Item(20, 18).Value = 1412.47028058455
$ws.Cells.Item(20, 19).Value = 0.001198404509125409
$ws.Cells.Item(20, 20).Value = 0.0005552553791775848

$ws.Cells.Item(21, 7).Value = 76.954105
$ws.Cells.Item(21, 8).Value = 153.90821
$ws.Cells.Item(21, 9).Value = 0.01070281688829022
$ws.Cells.Item(21, 10).Value = 0.007160758019481436
$ws.Cells.Item(21, 13).Value = 8.745164000000001
$ws.Cells.Item(21, 14).Value = 26.235492
$ws.Cells.Item(21, 15).Value = 0.2133957479992572
$ws.Cells.Item(21, 16).Value = 0.2216692531678338
$ws.Cells.Item(21, 17).Value = 672.97626869822
$ws.Cells.Item(21, 18).Value = 4037.85761218932
$ws.Cells.Item(21, 19).Value = 0.002283935615575774
$ws.Cells.Item(21, 20).Value = 0.001587319882294026

$ws.Cells.Item(22, 7).Value = 4655.195393666666
$ws.Cells.Item(22, 8).Value = 13965.586181
$ws.Cells.Item(22, 9).Value = 0.6474469929528313
$ws.Cells.Item(22, 10).Value = 0.6497650985763194
$ws.Cells.Item(22, 13).Value = 10.250695
$ws.Cells.Item(22, 14).Value = 30.752085
$ws.Cells.Item(22, 15).Value = 0.2501330709220828
$ws.Cells.Item(22, 16).Value = 0.259830908271274
$ws.Cells.Item(22, 17).Value = 47718.98814588191
$ws.Cells.Item(22, 18).Value = 429470.8933129373
$ws.Cells.Item(22, 19).Value = 0.1619479046065598
$ws.Cells.Item(22, 20).Value = 0.1688290557260589

$ws.Cells.Item(23, 7).Value = 4655.195393666666
$ws.Cells.Item(23, 8).Value = 13965.586181
$ws.Cells.Item(23, 9).Value = 0.6474469929528313
$ws.Cells.Item(23, 10).Value = 0.6497650985763194
$ws.Cells.Item(23, 15).Value = 0.2460808482884365
$ws.Cells.Item(23, 16).Value = 0.2556215780794031
$ws.Cells.Item(23, 17).Value = 46945.92777802806
$ws.Cells.Item(23, 18).Value = 422513.3500022526
$ws.Cells.Item(23, 19).Value = 0.1593243052476301
$ws.Cells.Item(23, 20).Value = 0.1660939798789977

$ws.Cells.Item(24, 7).Value = 4655.195393666666
$ws.Cells.Item(24, 8).Value = 13965.586181
$ws.Cells.Item(24, 9).Value = 0.6474469929528313
$ws.Cells.Item(24, 10).Value = 0.6497650985763194
$ws.Cells.Item(24, 13).Value = 7.311799000000001
$ws.Cells.Item(24, 14).Value = 21.935397
$ws.Cells.Item(24, 15).Value = 0.1784193889131434
$ws.Cells.Item(24, 16).Value = 0.1853368357235283
$ws.Cells.Item(24, 17).Value = 34037.85302421654
$ws.Cells.Item(24, 18).Value = 306340.6772179488
$ws.Cells.Item(24, 19).Value = 0.1155170968362964
$ws.Cells.Item(24, 20).Value = 0.1204254073337215

$ws.Cells.Item(25, 7).Value = 4655.195393666666
$ws.Cells.Item(25, 8).Value = 13965.586181
$ws.Cells.Item(25, 9).Value = 0.6474469929528313
$ws.Cells.Item(25, 10).Value = 0.6497650985763194
$ws.Cells.Item(25, 13).Value = 4.588677499999999
$ws.Cells.Item(25, 14).Value = 9.177354999999999
$ws.Cells.Item(25, 15).Value = 0.1119709438770801
$ws.Cells.Item(25, 16).Value = 0.07754142475796089
$ws.Cells.Item(25, 17).Value = 21361.19036102187
$ws.Cells.Item(25, 18).Value = 128167.1421661312
$ws.Cells.Item(25, 19).Value = 0.07249525091130574
$ws.Cells.Item(25, 20).Value = 0.05038371150160471

$ws.Cells.Item(26, 7).Value = 4655.195393666666
$ws.Cells.Item(26, 8).Value = 13965.586181
$ws.Cells.Item(26, 9).Value = 0.6474469929528313
$ws.Cells.Item(26, 10).Value = 0.6497650985763194
$ws.Cells.Item(26, 13).Value = 8.745164000000001
$ws.Cells.Item(26, 14).Value = 26.235492
$ws.Cells.Item(26, 15).Value = 0.2133957479992572
$ws.Cells.Item(26, 16).Value = 0.2216692531678338
$ws.Cells.Item(26, 17).Value = 20774.99288374904
$ws.Cells.Item(26, 18).Value = 186974.9359537413
$ws.Cells.Item(26, 19).Value = 0.0705058237080347
$ws.Cells.Item(26, 20).Value = 0.07350160947580005
